$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Working on the sheet BEFORE the row-13 deletion: reposition the data values
# that end up shifted by one row once the orphaned row 13 is removed.  Several
# of the target values already exist elsewhere on the sheet, so copy/paste the
# whole B:C cells (keeps the shared-string reuse + existing formatting intact)
# working from the bottom of each dependency chain upward so a source row is
# never clobbered before it has been read.
$ws.Range("B21:C21").Copy()
$ws.Range("B22:C22").PasteSpecial(-4163)

$ws.Range("B20:C20").Copy()
$ws.Range("B21:C21").PasteSpecial(-4163)

$ws.Range("B19:C19").Copy()
$ws.Range("B20:C20").PasteSpecial(-4163)

$ws.Range("B13:C13").Copy()
$ws.Range("B19:C19").PasteSpecial(-4163)

$ws.Range("B13:C13").Copy()
$ws.Range("B10:C10").PasteSpecial(-4163)

$ws.Range("B8:C8").Copy()
$ws.Range("B16:C16").PasteSpecial(-4163)

# "Semestral" is brand-new text that isn't present anywhere else on the
# sheet, so type it directly and then fix up the style (typing into a blank
# B/C cell otherwise inherits column A's bold style because of the
# overlapping <col> ranges in this workbook).
$ws.Range("B14").Value = "Semestral"
$ws.Range("B10").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C14").Value = "Semestral"
$ws.Range("C10").Copy()
$ws.Range("C14").PasteSpecial(-4122)

# Finally, remove the now-redundant orphan row (old row 13 held the "6376612
# - Daisy Rafaela da Silva" value with no label) which shifts everything
# below it up by one row.
$ws.Rows(13).Delete()
